$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.051202535629272
$ws.Range("B1").Value = 1.740820527076721
$ws.Range("C1").Value = 4.53245210647583
$ws.Range("D1").Value = 2.50859546661377
$ws.Range("E1").Value = 1.292465090751648
